# Implement and benchmark faster totient_sum
# - rename existing Python/PyPy version labels to include version numbers
# - add two new benchmark rows for the new totient_sum implementation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Python version label becomes explicit text "python 3.5"
# (was stored as the number 3.5; the column is text-formatted ("@"))
$ws.Range("C3").Value = "python 3.5"

# Row 4: PyPy version label gets an explicit version suffix
$ws.Range("C4").Value = "pypy 5.1.2"

# New row 5: totient_sum (faster implementation) timings on python 3.5
$ws.Range("A5").Value = "totient_sum "
$ws.Range("B5").Value = "laptop"
$ws.Range("C5").Value = "python 3.5"
$ws.Range("G5").Value = "0.018"
$ws.Range("H5").Value = "0.092"
$ws.Range("I5").Value = "0.55"
$ws.Range("J5").Value = "2.35"
$ws.Range("K5").Value = "11.47"
$ws.Range("L5").Value = "57.91"

# New row 6: totient_sum (faster implementation) timings on pypy 5.1.2
$ws.Range("A6").Value = "totient_sum "
$ws.Range("B6").Value = "laptop"
$ws.Range("C6").Value = "pypy 5.1.2"
$ws.Range("G6").Value = "0.016"
$ws.Range("H6").Value = "0.029"
$ws.Range("I6").Value = "0.065"
$ws.Range("J6").Value = "0.21"
$ws.Range("K6").Value = "0.99"
$ws.Range("L6").Value = "4.61"
$ws.Range("M6").Value = "23.39"

# Update the active selection to reflect where the author was working
[void]$ws.Range("D7").Select()
